$d = $word.ActiveDocument

# Update the submission timestamp (date unchanged, time corrected).
$d.Content.Find.Execute(
    "10/3/2013, 01:02:33", $true, $false, $false, $false, $false,
    $true, 1, $false, "10/3/2013, 17:50:28", 2
) | Out-Null

# Update the reported file size.
$d.Content.Find.Execute(
    "2,733", $true, $false, $false, $false, $false,
    $true, 1, $false, "2,729", 2
) | Out-Null

# Update the confirmation code (correct project id suffix).
$d.Content.Find.Execute(
    "236321-303266019-6584280", $true, $false, $false, $false, $false,
    $true, 1, $false, "236321-303266019-7677834", 2
) | Out-Null
